# Update "Generate Report for Handback" timestamps.
# These cells hold plain text timestamps (displayed with a date/time
# number format), recording when each handback XLIFF/report was (re)generated.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 18d01c8d... row
$overview.Range("G3").Value = "2016-08-15 22:43:03"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$zhcn.Range("H3").Value = "2016-08-15 22:42:56"
$zhcn.Range("K3").Value = "2016-08-15 22:43:28"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$dede.Range("H3").Value = "2016-08-15 22:43:03"
$dede.Range("K3").Value = "2016-08-15 22:43:35"
